$d = $word.ActiveDocument

# Map of old bookmark name -> new (anonymized) bookmark name, as required by the
# commit "added guiding questions to school systems course": the heading
# bookmarks were re-keyed to opaque hash-based names while keeping their
# heading text, style and position untouched.
$renames = @(
    @{ Old = "psycho-sociological-aspects-of-social-media";            New = "Xbcf7ca2f8c0ac612646ad1233e97e6bf5068f16" },
    @{ Old = "social-media-privacy-safety-and-self-presentation";      New = "Xe93da2363c2ec44e8b9a0af3a23709bf30d345f" },
    @{ Old = "social-aspects-of-virtual-simulations-and-games";        New = "Xe8654de4d34414e241b6e84df987c969bd6f563" },
    @{ Old = "the-video-revolution-and-the-power-of-video";            New = "Xd24ca5b46d5a47dd66068533f3f46eedd089483" },
    @{ Old = "critical-perspectives-on-social-media---case-studies";   New = "Xf3a7c7157b99f7dfaab7427ca23cd16b8f5f256" },
    @{ Old = "summaries-and-evaluations-of-social-media-in-education"; New = "X374b11a8d083674120493be9d74371d5b995bb2" }
)

function Escape-Xml($text) {
    $text = $text -replace "&", "&amp;"
    $text = $text -replace "<", "&lt;"
    $text = $text -replace ">", "&gt;"
    return $text
}

foreach ($pair in $renames) {
    $oldName = $pair.Old
    $newName = $pair.New

    $bm = $d.Bookmarks($oldName)
    $bmRange = $bm.Range
    $headingText = $bmRange.Text

    $para = $bmRange.Paragraphs(1)
    $pRange = $para.Range

    $escapedText = Escape-Xml($headingText)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
           '<w:bookmarkStart w:id="0" w:name="' + $newName + '"/>' +
           '<w:r><w:t xml:space="preserve">' + $escapedText + '</w:t></w:r>' +
           '<w:bookmarkEnd w:id="0"/>' +
           '</w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $pRange.InsertXML($xml)
}

Write-Host "Renamed" $renames.Count "bookmarks"
